$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: date 44995 -> 45008, quality Primera -> Especial, values updated
$ws.Range("D2").Value = 45008
$ws.Range("L2").Value = "Especial"
$ws.Range("M2").Value = 60
$ws.Range("N2").Value = 7000
$ws.Range("O2").Value = 7000
$ws.Range("P2").Value = 7000
$ws.Range("S2").Value = 3500

# Row 3: date stays 45008, quality Especial -> Primera, values updated
$ws.Range("L3").Value = "Primera"
$ws.Range("N3").Value = 6000
$ws.Range("O3").Value = 6000
$ws.Range("P3").Value = 6000
$ws.Range("S3").Value = 3000

# Row 4: date 45008 -> 44995, quality stays Primera, values updated
$ws.Range("D4").Value = 44995
$ws.Range("M4").Value = 100
$ws.Range("N4").Value = 5500
$ws.Range("P4").Value = 5750
$ws.Range("S4").Value = 2875
